# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect the latest generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - row => new F value
$wsExhibition = $wb.Worksheets.Item("展览")
$exhibitionUpdates = @{
    2  = 305
    3  = 67
    4  = 3705
    6  = 441
    9  = 180
    10 = 99
    11 = 80
    12 = 1369
    13 = 242
    14 = 2148
    15 = 157
}
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Sheet "全部类型" (all types) - row => new F value
$wsAll = $wb.Worksheets.Item("全部类型")
$allUpdates = @{
    2  = 305
    3  = 67
    4  = 3705
    6  = 441
    10 = 180
    11 = 99
    12 = 80
    15 = 1369
    16 = 242
    17 = 2148
    18 = 157
}
foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}
